# "merge fail - copy in files"
# Adds a new "POX" (hydrocarbon partial oxidation) worksheet between "IEA Data"
# and "HPEbP", wires the HPEbP "thermochemical water splitting" row over to the
# new "hydrocarbon partial oxidation" pathway sourced from the POX sheet, and
# replaces the "Thermochemical Water Splitting" write-up on the About sheet
# with a new "Efficiency Data (POX)" source block.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "POX" worksheet right after "IEA Data" (before "HPEbP").
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("IEA Data")
$pox = $wb.Worksheets.Add($null, $afterSheet)
$pox.Name = "POX"

$pox.Range("B2").Value = "From El-Shafie et al."
$pox.Range("B3").Value = "60-75%"
$pox.Range("C3").Value = "efficiency"
$pox.Range("B5").Formula = "=AVERAGE(0.6,0.75)"

$pox.Range("B6").Select()

# ---------------------------------------------------------------------------
# 2. HPEbP: swap the "thermochemical water splitting" row for
#    "hydrocarbon partial oxidation", sourced off the new POX sheet.
# ---------------------------------------------------------------------------
$hpebp = $wb.Worksheets.Item("HPEbP")

$hpebp.Range("A6").Value = "hydrocarbon partial oxidation"
$hpebp.Range("B6").Formula = "=POX!B5"

$hpebp.Range("B6").Select()

# ---------------------------------------------------------------------------
# 3. About sheet: remove the old "Thermochemical Water Splitting" note block
#    and add a new "Efficiency Data (POX)" source block; the "Notes" block
#    moves down to make room.
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Wipe out the old rows 19-30 (Notes + Thermochemical Water Splitting block);
# they get rebuilt below at their new positions.
$about.Range("A19:A30").Clear()

# New "Efficiency Data (POX)" source block (rows 17-23), matching the style
# of the existing "Sources" (row 3) / "Efficiency Data" (row 10) blocks.
$about.Range("B17").Value = "Efficiency Data (POX)"
$about.Range("B10").Copy()
$about.Range("B17").PasteSpecial(-4122)

$about.Range("B18").Value = "Journal of Power and Energy Engineering"

$about.Range("B19").Value = 2019
$about.Range("B12").Copy()
$about.Range("B19").PasteSpecial(-4122)

$about.Range("B20").Value = "Hydrogen Production Technologies Overview"
$about.Range("B21").Value = "Mostafa El-Shafie et al."

$about.Range("B22").Value = "https://www.scirp.org/journal/paperinformation?paperid=90227"
$about.Range("B14").Copy()
$about.Range("B22").PasteSpecial(-4122)

$about.Range("B23").Value = "Section 3.2"

# "Notes" block moves from rows 19-21 down to rows 26-28 (style copied from
# A3, which carries the same bold "section header" style as the old A19).
$about.Range("A26").Value = "Notes"
$about.Range("A3").Copy()
$about.Range("A26").PasteSpecial(-4122)

$about.Range("A27").Value = "This variable expresses the amount of energy input of each"
$about.Range("A28").Value = "source fuel to produce one unit of energy of hydrogen."

# A leftover styled-but-empty cell remains at A30 (same bold style as the
# section headers).
$about.Range("A3").Copy()
$about.Range("A30").PasteSpecial(-4122)
$about.Range("A30").ClearContents()

$about.Range("B26").Select()

# ---------------------------------------------------------------------------
# 4. Restore view state: HPEbP becomes the active/selected tab.
# ---------------------------------------------------------------------------
$about.Activate()
$about.Range("B26").Select()

$iea = $wb.Worksheets.Item("IEA Data")
$iea.Activate()
$iea.Range("I6").Select()

$pox.Activate()
$pox.Range("B6").Select()

$hpebp.Activate()
$hpebp.Range("B6").Select()
